$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster FAPs -> M1, recompute stats ---
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.822099333333334
$ws.Range("H2").Value = 8.466298
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1556706666666667
$ws.Range("N2").Value = 0.467012
$ws.Range("O2").Value = 0.00151537522169743
$ws.Range("P2").Value = 0.00151537522169743
$ws.Range("Q2").Value = 0.4393180846195556
$ws.Range("R2").Value = 3.953862761576
$ws.Range("S2").Value = 0.00151537522169743
$ws.Range("T2").Value = 0.00151537522169743

# --- Row 3: Target cluster M1 -> FAPs, recompute stats ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.822099333333334
$ws.Range("H3").Value = 8.466298
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 37.12743
$ws.Range("N3").Value = 111.38229
$ws.Range("O3").Value = 0.3614167567469731
$ws.Range("P3").Value = 0.3614167567469732
$ws.Range("Q3").Value = 104.77729545138
$ws.Range("R3").Value = 942.99565906242
$ws.Range("S3").Value = 0.3614167567469731
$ws.Range("T3").Value = 0.3614167567469732

# --- Row 4: Target cluster M2, recompute stats ---
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.822099333333334
$ws.Range("H4").Value = 8.466298
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 17.74257933333334
$ws.Range("N4").Value = 53.227738
$ws.Range("O4").Value = 0.1727150378838289
$ws.Range("P4").Value = 0.1727150378838289
$ws.Range("Q4").Value = 50.07132130821379
$ws.Range("R4").Value = 450.641891773924
$ws.Range("S4").Value = 0.1727150378838289
$ws.Range("T4").Value = 0.1727150378838289

# --- Row 5: Target cluster Neutro, recompute stats ---
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.822099333333334
$ws.Range("H5").Value = 8.466298
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 18.87700766666667
$ws.Range("N5").Value = 56.631023
$ws.Range("O5").Value = 0.1837581240601467
$ws.Range("P5").Value = 0.1837581240601467
$ws.Range("Q5").Value = 53.27279075142823
$ws.Range("R5").Value = 479.455116762854
$ws.Range("S5").Value = 0.1837581240601467
$ws.Range("T5").Value = 0.1837581240601467

# --- Row 6: Target cluster sCs, recompute stats ---
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.822099333333334
$ws.Range("H6").Value = 8.466298
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.878365
$ws.Range("N6").Value = 14.635095
$ws.Range("O6").Value = 0.04748841642225027
$ws.Range("P6").Value = 0.04748841642225028
$ws.Range("Q6").Value = 13.76723061425667
$ws.Range("R6").Value = 123.90507552831
$ws.Range("S6").Value = 0.04748841642225027
$ws.Range("T6").Value = 0.04748841642225028

# --- Row 7 (new): Target cluster ECs ---
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Sdc4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.822099333333334
$ws.Range("H7").Value = 8.466298
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.94642
$ws.Range("N7").Value = 71.83926
$ws.Range("O7").Value = 0.2331062896651035
$ws.Range("P7").Value = 0.2331062896651035
$ws.Range("Q7").Value = 67.57917591772001
$ws.Range("R7").Value = 608.21258325948
$ws.Range("S7").Value = 0.2331062896651035
$ws.Range("T7").Value = 0.2331062896651035
